# Append a new time-tracking entry (row 9) to the sheet and update the
# active selection, matching the "Card image popper is now faster and more
# convenient" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data: Date (A9), Task description (B9), Time in hours (C9)
$ws.Range("A9").Value = "28.06.2019"
$ws.Range("B9").Value = "Card image Popper now opens up nicely above the card text `nline when hovering the mouse over it"
$ws.Range("C9").Value = 3

# The task-description column uses the wrapped-text style seen on the other
# entries (e.g. B5, B7, B8) with a taller (30pt) row to fit the two lines.
$ws.Range("B9").WrapText = $true
$ws.Rows.Item(9).RowHeight = 30

# Move the active selection to match the post-edit workbook state.
$ws.Range("E7").Select()
